$wb = $excel.ActiveWorkbook

# --- Add the two new worksheets after the last existing sheet ("promotion") ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheetA = $wb.Worksheets.Add($null, $lastSheet)
$sheetA.Name = "mergeA"
$sheetB = $wb.Worksheets.Add($null, $sheetA)
$sheetB.Name = "mergeB"

# --- mergeA data ---
$dataA = @(
    @("Key", "Address.State", "Address.City", "Address.TEL(Int)"),
    @(1, "Some", "New", "555;1111;2222"),
    @(2, "Place", "York", "555;3333;4444"),
    @(3, "Beyond", "Los", "555;5555;6666"),
    @(4, "Rainbow", "Angeles", "555;7777;8888")
)
for ($r = 0; $r -lt $dataA.Length; $r++) {
    for ($c = 0; $c -lt $dataA[$r].Length; $c++) {
        $sheetA.Cells.Item($r + 1, $c + 1).Value = $dataA[$r][$c]
    }
}
$sheetA.Columns.Item(2).ColumnWidth = 9.2

# --- mergeB data ---
$dataB = @(
    @("Key", "Name", "Property.[1,A]", "Property.[2,A]"),
    @(1, "My", "Out", "think"),
    @(2, "Name", "of", "of"),
    @(3, "Is", "idea", "anything"),
    @(4, "Jonas", "cannot", "funny")
)
for ($r = 0; $r -lt $dataB.Length; $r++) {
    for ($c = 0; $c -lt $dataB[$r].Length; $c++) {
        $sheetB.Cells.Item($r + 1, $c + 1).Value = $dataB[$r][$c]
    }
}

# --- Selections on each sheet ---
$sheetB.Activate() | Out-Null
$sheetB.Range("D30").Select() | Out-Null

$sheetA.Activate() | Out-Null
$sheetA.Range("D2").Select() | Out-Null
